$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert 3 new columns at M:O, pushing the existing "% of Total Deaths" column from M to P
$ws.Columns("M:O").Insert()

# Set the 3 new header cells. Setting order (O, N, then M) reproduces the original
# shared-string insertion order: Standard Error, Upper CI, Lower CI.
$ws.Cells.Item(1, 15).Value2 = "Age Adjusted Rate Standard Error"
$ws.Cells.Item(1, 14).Value2 = "Age Adjusted Rate Upper 95% Confidence Interval"
$ws.Cells.Item(1, 13).Value2 = "Age Adjusted Rate Lower 95% Confidence Interval"

# Data for new columns M (Lower 95% CI), N (Upper 95% CI), O (Standard Error),
# one row per array entry for rows 2..100
$newData = @(
    @(0.06167, 0.10068, 0.00961),
    @(0.0794, 0.12737, 0.01183),
    @(0.05921, 0.10179, 0.01047),
    @(0.06114, 0.12936, 0.01658),
    @(0.07057, 0.14565, 0.01814),
    @(0.07231, 0.16221, 0.02169),
    @(0.08422, 0.10631, 0.00564),
    @(0.08023, 0.12487, 0.01102),
    @(0.08204, 0.1325, 0.01249),
    @(0.07776, 0.12603, 0.01194),
    @(0.07492, 0.15285, 0.01902),
    @(0.05345, 0.12192, 0.01662),
    @(0.0792, 0.16984, 0.02202),
    @(0.09066, 0.11356, 0.00584),
    @(0.07159, 0.11237, 0.01009),
    @(0.0791, 0.12605, 0.01163),
    @(0.07081, 0.11892, 0.01188),
    @(0.05926, 0.12709, 0.01639),
    @(0.07844, 0.15356, 0.01835),
    @(0.04567, 0.11547, 0.01653),
    @(0.07941, 0.09954, 0.00514),
    @(0.08087, 0.12303, 0.01047),
    @(0.08654, 0.13294, 0.01153),
    @(0.08015, 0.129, 0.01203),
    @(0.08409, 0.16166, 0.01902),
    @(0.07687, 0.15866, 0.02001),
    @(0.04223, 0.10679, 0.01527),
    @(0.08921, 0.11029, 0.00538),
    @(0.07127, 0.11033, 0.00966),
    @(0.07663, 0.12028, 0.01079),
    @(0.10687, 0.1598, 0.0135),
    @(0.08073, 0.16471, 0.02044),
    @(0.08839, 0.18033, 0.02236),
    @(0.06299, 0.13908, 0.01852),
    @(0.09589, 0.11795, 0.00563),
    @(0.07013, 0.10885, 0.00962),
    @(0.0827, 0.12736, 0.01105),
    @(0.06565, 0.10758, 0.01032),
    @(0.06021, 0.13096, 0.01723),
    @(0.04564, 0.10411, 0.01401),
    @(0.04322, 0.11209, 0.01628),
    @(0.08197, 0.10235, 0.0052),
    @(0.07358, 0.11074, 0.00948),
    @(0.0924, 0.13791, 0.01161),
    @(0.09083, 0.14215, 0.0127),
    @(0.14009, 0.23849, 0.02431),
    @(0.0433, 0.09712, 0.013),
    @(0.08859, 0.17512, 0.021),
    @(0.09566, 0.11671, 0.00537),
    @(0.08304, 0.12266, 0.01011),
    @(0.08735, 0.1296, 0.01078),
    @(0.08055, 0.12537, 0.01108),
    @(0.07554, 0.14932, 0.01797),
    @(0.07111, 0.14349, 0.01772),
    @(0.0561, 0.13278, 0.01858),
    @(0.09572, 0.11698, 0.00542),
    @(0.0887, 0.12769, 0.00995),
    @(0.09477, 0.13814, 0.01106),
    @(0.07875, 0.1229, 0.01094),
    @(0.09479, 0.17045, 0.01862),
    @(0.07788, 0.14844, 0.01718),
    @(0.05917, 0.13272, 0.01783),
    @(0.10086, 0.12195, 0.00538),
    @(0.08049, 0.11704, 0.00932),
    @(0.12137, 0.16978, 0.01235),
    @(0.08235, 0.12651, 0.01096),
    @(0.05929, 0.12096, 0.0151),
    @(0.08952, 0.16423, 0.01824),
    @(0.05624, 0.12615, 0.01699),
    @(0.10245, 0.12347, 0.00536),
    @(0.08621, 0.12649, 0.01028),
    @(0.10113, 0.14534, 0.01128),
    @(0.08961, 0.13482, 0.01153),
    @(0.07275, 0.14679, 0.01803),
    @(0.07262, 0.14501, 0.01756),
    @(0.05975, 0.1363, 0.01837),
    @(0.09624, 0.11646, 0.00516),
    @(0.08374, 0.12163, 0.00967),
    @(0.08034, 0.11861, 0.00976),
    @(0.08821, 0.13388, 0.01136),
    @(0.07637, 0.14951, 0.01792),
    @(0.06561, 0.12969, 0.01568),
    @(0.06247, 0.14506, 0.0198),
    @(0.0946, 0.115, 0.00521),
    @(0.10057, 0.14131, 0.01039),
    @(0.12439, 0.17033, 0.01172),
    @(0.06879, 0.11072, 0.01035),
    @(0.09946, 0.17563, 0.01866),
    @(0.0707, 0.14424, 0.01794),
    @(0.07671, 0.15831, 0.01992),
    @(0.10418, 0.1245, 0.00518),
    @(0.06911, 0.10201, 0.00839),
    @(0.09765, 0.13999, 0.0108),
    @(0.08968, 0.13428, 0.01138),
    @(0.08911, 0.16126, 0.01771),
    @(0.09491, 0.16573, 0.01743),
    @(0.04905, 0.1161, 0.01597),
    @(0.10564, 0.12728, 0.00552),
    @(0.102, 0.1076, 0.00143)
)

for ($i = 0; $i -lt $newData.Length; $i++) {
    $r = $i + 2
    $vals = $newData[$i]
    $ws.Cells.Item($r, 13).Value2 = $vals[0]
    $ws.Cells.Item($r, 14).Value2 = $vals[1]
    $ws.Cells.Item($r, 15).Value2 = $vals[2]
}

Write-Host "done"